$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "coal_coke"
$ws.Range("B6").Value = "coal_coke"

$ws.Range("A7").Select()
